# Insert two new list-item paragraphs after the
# "Verallgemeinerung auf Anzahl Typen > 2" bullet in the "Code" section:
#   1) a red "Untersuche Zusammenhang Temperatur – Nbar" bullet
#   2) a new empty bullet (same list) directly below it

$d = $word.ActiveDocument

# Locate the anchor paragraph via Find so we are not dependent on
# hard-coded character offsets.
$anchor = $d.Content
$found = $anchor.Find.Execute("Verallgemeinerung auf Anzahl Typen > 2", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Anchor paragraph 'Verallgemeinerung auf Anzahl Typen > 2' not found"
}

# $anchor now spans just the found text; move one character past its end
# so the insertion point sits right after the paragraph mark that closes
# that paragraph (i.e. the very start of the following paragraph).
$insPos = $anchor.End + 1
$insertionPoint = $d.Range($insPos, $insPos)

$newParagraphsXml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:rPr><w:color w:val="FF0000"/><w:lang w:val="de-DE"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="FF0000"/><w:lang w:val="de-DE"/></w:rPr><w:t xml:space="preserve">Untersuche Zusammenhang Temperatur – </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:color w:val="FF0000"/><w:lang w:val="de-DE"/></w:rPr><w:t>Nbar</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:rPr><w:lang w:val="de-DE"/></w:rPr></w:pPr></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$insertionPoint.InsertXML($newParagraphsXml)
